$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $result = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                       $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $findText"
    }
}

# --- "Programa resumido" paragraph: split the run at each numbered item boundary ---
Replace-Text "Produção2. Visão estratégica da Produção." "Produção^l2. Visão estratégica da Produção."
Replace-Text "Produção.3. Projeto em Gestão da Produção." "Produção.^l3. Projeto em Gestão da Produção."
Replace-Text "Produção.4. Planejamento e Controle da Produção" "Produção.^l4. Planejamento e Controle da Produção"

# --- "Programa" paragraph: split the run at each numbered item boundary ---
Replace-Text "Produção2 - Visão estratégica de produção:" "Produção^l2 - Visão estratégica de produção:"
Replace-Text "Produto/Serviço.3 – Projeto em Gestão da Produção:" "Produto/Serviço.^l3 – Projeto em Gestão da Produção:"
Replace-Text "Arranjo Físico.4 - Planejamento e Controle da Produção:" "Arranjo Físico.^l4 - Planejamento e Controle da Produção:"

# --- "Bibliografia" paragraph: split the run after each reference ---
Replace-Text "Atlas, 2009.CHASE, R. B. E JACOBS, F.R." "Atlas, 2009.^lCHASE, R. B. E JACOBS, F.R."
Replace-Text "Bookman. 2009.CORREA, H.L.; CORREA, C.A." "Bookman. 2009.^lCORREA, H.L.; CORREA, C.A."
